# Integración de BD de Finnegans de productos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Unidades" shifts from B to C)
$ws.Columns.Item(2).Insert()

# New header for the inserted column, matching the formatting of the other header cells
$ws.Range("B1").Value = "Nombre del Producto"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4160
$ws.Range("B1").Borders.LineStyle = 1

# Replace the old data rows (previously rows 2-3) with the new product rows (2-6)
$data = @(
    @("ROPA001", "BOTA PETROLERA-T 40-MARRON-MASC-BORIS C/CORDONES 3703 -21%", 2),
    @("ROPA007", "BOTA PETROLERA-T 38-MARRON-MASC-BORIS LOW LINE 3910X/3910D-21%", 1),
    @("ROPA009", "BOTA PETROLERA-T 40-MARRON-MASC-BORIS LOW LINE 3910X/3910D-21%", 1),
    @("ROPA123", "ZAPATO -T 38-MARRON-MASC-BORIS 3161D MD-21%", 1),
    @("ROPA112", "BOTIN -T 42-MARRON-MASC-FUNCIONAL TERRA BROWN-21%", 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
